# "upgrading to the latest SHAFT_Engine version"
#
# The "Expected First Search Result Link" cell (B4 on the SimpleSearch
# sheet) used to hold the full github URL as both its displayed text and
# as an external hyperlink. The new version just shows the short repo
# name "SHAFT_ENGINE" as plain text with the hyperlink removed (and the
# leftover "hyperlink blue" font formatting cleared back to the normal
# black font).
#
# Along with that, the sheet's active selection moved from B2 to A6, and
# row 4's height settled back down to the sheet's normal 13.8pt (it used
# to be taller to fit the long URL).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimpleSearch")

# Drop the external hyperlink that used to live on B4.
foreach ($hl in $ws.Hyperlinks) {
    $hl.Delete()
}

# Replace the long repo URL with just the short project name.
$cell = $ws.Range("B4")
$cell.Value = "SHAFT_ENGINE"

# The hyperlink styling (blue font) is no longer appropriate now that the
# cell is plain text - put the font color back to automatic/black.
$cell.Font.Color = 0

# Row 4 no longer needs the extra height that the long URL required.
$ws.Rows.Item(4).RowHeight = 13.8

# Active cell/selection moved to A6.
$ws.Range("A6").Select()
